$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.729.02"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "1.805.91"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'232.52"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'0.5925"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'0.2775"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "'0.06828"
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").Value = "'0.07512"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").Value = "1.806.63"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "'4.753"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "'0.6227"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "2.051.70"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "'0.000009276"
$ws.Range("E16").Value = "  -6.80%  "
$ws.Range("D17").Value = "'75.61"
$ws.Range("E17").Value = "  -4.12%  "
$ws.Range("D18").Value = "28.699.01"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").Value = "'5.472"
$ws.Range("E19").Value = "  -6.58%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "'210.96"
$ws.Range("E21").Value = "  -6.91%  "
$ws.Range("D22").Value = "'11.49"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "'6.824"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "'153.97"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "'7.870"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").Value = "'16.44"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").Value = "'1.437"
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("D30").Value = "'0.06156"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").Value = "'1.428"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "'3.783"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").Value = "'3.748"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "'1.061"
$ws.Range("E35").Value = "  -5.43%  "
$ws.Range("D36").Value = "'0.6421"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'2.497"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("D38").Value = "'2.720"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").Value = "'6.573"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").Value = "1.145.28"
$ws.Range("E41").Value = "  -5.62%  "
$ws.Range("D42").Value = "'0.8830"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").Value = "'1.007"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "'100.05"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").Value = "1.959.48"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").Value = "'60.48"
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("D47").Value = "'0.00000000112"
$ws.Range("E47").Value = "  -5.18%  "
$ws.Range("D48").Value = "'1.596"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'8.352"
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").Value = "'0.05466"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'0.4479"
$ws.Range("E51").Value = "  -1.70%  "
